# Update the pl_mw results table on Sheet1 with the recalculated values
# from the "case with 380 kV" run (rows 2-25, columns B-O excluding the
# always-zero F/I/L/N columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3391143333193156
$ws.Range("C2").Value = 0.03373466785679113
$ws.Range("D2").Value = 0.6371284573104248
$ws.Range("E2").Value = 0.2525459256420177
$ws.Range("G2").Value = 0.9775633400828809
$ws.Range("H2").Value = 1.016628841809762
$ws.Range("J2").Value = 0.1260186915024093
$ws.Range("K2").Value = 0.2957287194307412
$ws.Range("M2").Value = 0.2675125823171811
$ws.Range("O2").Value = 4.03728545767737

$ws.Range("B3").Value = 0.3066594699521659
$ws.Range("C3").Value = 0.03034087177651656
$ws.Range("D3").Value = 0.6317625143170176
$ws.Range("E3").Value = 0.2512538757166318
$ws.Range("G3").Value = 0.9820653683865075
$ws.Range("H3").Value = 1.022953157778204
$ws.Range("J3").Value = 0.1259393466668008
$ws.Range("K3").Value = 0.2625023828075257
$ws.Range("M3").Value = 0.2561742159976461
$ws.Range("O3").Value = 4.059775269391679

$ws.Range("B4").Value = 0.2867846183939378
$ws.Range("C4").Value = 0.02824515716683607
$ws.Range("D4").Value = 0.6287632993262093
$ws.Range("E4").Value = 0.2505758187692173
$ws.Range("G4").Value = 0.9853685517947781
$ws.Range("H4").Value = 1.027230046640092
$ws.Range("J4").Value = 0.125945988430054
$ws.Range("K4").Value = 0.2421010022130332
$ws.Range("M4").Value = 0.2493211293329551
$ws.Range("O4").Value = 4.075540495075415

$ws.Range("B5").Value = 0.2786991507011294
$ws.Range("C5").Value = 0.02738819561297845
$ws.Range("D5").Value = 0.6276155524782325
$ws.Range("E5").Value = 0.2503285365435666
$ws.Range("G5").Value = 0.986850118470592
$ws.Range("H5").Value = 1.029072003459717
$ws.Range("J5").Value = 0.1259626297231016
$ws.Range("K5").Value = 0.2337876871110325
$ws.Range("M5").Value = 0.2465559562180957
$ws.Range("O5").Value = 4.082457019176417

$ws.Range("B6").Value = 0.2773574062314594
$ws.Range("C6").Value = 0.0272457217776747
$ws.Range("D6").Value = 0.6274294715472593
$ws.Range("E6").Value = 0.2502892303244764
$ws.Range("G6").Value = 0.9871043142371221
$ws.Range("H6").Value = 1.029383846201938
$ws.Range("J6").Value = 0.1259662350346353
$ws.Range("K6").Value = 0.2324073066684491
$ws.Range("M6").Value = 0.2460984684326348
$ws.Range("O6").Value = 4.08363522580089

$ws.Range("B7").Value = 0.286675518861415
$ws.Range("C7").Value = 0.02823361171915906
$ws.Range("D7").Value = 0.62874751873116
$ws.Range("E7").Value = 0.2505723662249935
$ws.Range("G7").Value = 0.9853879841485025
$ws.Range("H7").Value = 1.027254486594735
$ws.Range("J7").Value = 0.1259461564204862
$ws.Range("K7").Value = 0.2419888836242734
$ws.Range("M7").Value = 0.2492837255963067
$ws.Range("O7").Value = 4.075631781229944

$ws.Range("B8").Value = 0.3279132622310783
$ws.Range("C8").Value = 0.03256699024005627
$ws.Range("D8").Value = 0.6352170385336251
$ws.Range("E8").Value = 0.2520765301814869
$ws.Range("G8").Value = 0.9790037713902251
$ws.Range("H8").Value = 1.018727795777636
$ws.Range("J8").Value = 0.1259798515795367
$ws.Range("K8").Value = 0.2842726143876746
$ws.Range("M8").Value = 0.2635806521117559
$ws.Range("O8").Value = 4.044633962246678

$ws.Range("B9").Value = 0.4091805791561001
$ws.Range("C9").Value = 0.04096825513074975
$ws.Range("D9").Value = 0.6502430090611995
$ws.Range("E9").Value = 0.2559391629600896
$ws.Range("G9").Value = 0.9707622152117352
$ws.Range("H9").Value = 1.005127635586518
$ws.Range("J9").Value = 0.1264847092071193
$ws.Range("K9").Value = 0.3671716126931415
$ws.Range("M9").Value = 0.2924733319644162
$ws.Range("O9").Value = 3.999367302171464

$ws.Range("B10").Value = 0.4691150880456405
$ws.Range("C10").Value = 0.04707978096541865
$ws.Range("D10").Value = 0.6627034510894418
$ws.Range("E10").Value = 0.259332072155253
$ws.Range("G10").Value = 0.9673185701205966
$ws.Range("H10").Value = 0.9970336727870972
$ws.Range("J10").Value = 0.1271226965732666
$ws.Range("K10").Value = 0.4280488609504118
$ws.Range("M10").Value = 0.314217047532189
$ws.Range("O10").Value = 3.975570074470028

$ws.Range("B11").Value = 0.4964270995037907
$ws.Range("C11").Value = 0.04984645405039601
$ws.Range("D11").Value = 0.6686795764955207
$ws.Range("E11").Value = 0.2609958297924848
$ws.Range("G11").Value = 0.9663198455479574
$ws.Range("H11").Value = 0.993762779671556
$ws.Range("J11").Value = 0.127470849863542
$ws.Range("K11").Value = 0.455734085148265
$ws.Range("M11").Value = 0.3242198407195005
$ws.Range("O11").Value = 3.966798289401652

$ws.Range("B12").Value = 0.5067758799249305
$ws.Range("C12").Value = 0.05089213521196712
$ws.Range("D12").Value = 0.6709867183712674
$ws.Range("E12").Value = 0.2616431146190052
$ws.Range("G12").Value = 0.9660233555568283
$ws.Range("H12").Value = 0.9925832200587905
$ws.Range("J12").Value = 0.1276110064419882
$ws.Range("K12").Value = 0.466216185627502
$ws.Range("M12").Value = 0.3280235219737904
$ws.Range("O12").Value = 3.963771926275683

$ws.Range("B13").Value = 0.5045468122692682
$ws.Range("C13").Value = 0.05066701883909275
$ws.Range("D13").Value = 0.6704878745389919
$ws.Range("E13").Value = 0.2615029434894893
$ws.Range("G13").Value = 0.9660835752387555
$ws.Range("H13").Value = 0.9928346336825626
$ws.Range("J13").Value = 0.1275804514345822
$ws.Range("K13").Value = 0.4639587613373237
$ws.Range("M13").Value = 0.3272036296431793
$ws.Range("O13").Value = 3.964410573176593

$ws.Range("B14").Value = 0.4972783763003292
$ws.Range("C14").Value = 0.04993252324061359
$ws.Range("D14").Value = 0.6688685032147816
$ws.Range("E14").Value = 0.2610487367224579
$ws.Range("G14").Value = 0.9662938153262814
$ws.Range("H14").Value = 0.9936645533387605
$ws.Range("J14").Value = 0.1274822139820841
$ws.Range("K14").Value = 0.456596490780413
$ws.Range("M14").Value = 0.3245324554753282
$ws.Range("O14").Value = 3.966543389987947

$ws.Range("B15").Value = 0.4928270564192019
$ws.Range("C15").Value = 0.04948236157989072
$ws.Range("D15").Value = 0.6678823318805485
$ws.Range("E15").Value = 0.2607727678865928
$ws.Range("G15").Value = 0.9664332355300189
$ws.Range("H15").Value = 0.9941805922998981
$ws.Range("J15").Value = 0.1274231236695869
$ws.Range("K15").Value = 0.4520866515242403
$ws.Range("M15").Value = 0.3228983414906494
$ws.Range("O15").Value = 3.967888261525985

$ws.Range("B16").Value = 0.4673310955655836
$ws.Range("C16").Value = 0.04689869594457718
$ws.Range("D16").Value = 0.6623190783432165
$ws.Range("E16").Value = 0.2592257587921765
$ws.Range("G16").Value = 0.9673952695941921
$ws.Range("H16").Value = 0.99725570060005
$ws.Range("J16").Value = 0.1271011086263201
$ws.Range("K16").Value = 0.4262393579176091
$ws.Range("M16").Value = 0.3135655683833463
$ws.Range("O16").Value = 3.976184657068956

$ws.Range("B17").Value = 0.4517019481013165
$ws.Range("C17").Value = 0.04531020411840814
$ws.Range("D17").Value = 0.6589849331427899
$ws.Range("E17").Value = 0.2583075008604681
$ws.Range("G17").Value = 0.9681309148895707
$ws.Range("H17").Value = 0.9992474289944937
$ws.Range("J17").Value = 0.1269183916510528
$ws.Range("K17").Value = 0.4103804230065009
$ws.Range("M17").Value = 0.3078686298631936
$ws.Range("O17").Value = 3.981800229312171

$ws.Range("B18").Value = 0.4427169778172981
$ws.Range("C18").Value = 0.04439527926682274
$ws.Range("D18").Value = 0.657096202290262
$ws.Range("E18").Value = 0.2577906696093564
$ws.Range("G18").Value = 0.9686074799953417
$ws.Range("H18").Value = 1.000431715754829
$ws.Range("J18").Value = 0.1268187509924559
$ws.Range("K18").Value = 0.4012580478966754
$ws.Range("M18").Value = 0.3046024046405194
$ws.Range("O18").Value = 3.985223461088964

$ws.Range("B19").Value = 0.4396756102358097
$ws.Range("C19").Value = 0.04408528570517944
$ws.Range("D19").Value = 0.6564616925917051
$ws.Range("E19").Value = 0.2576176260912817
$ws.Range("G19").Value = 0.968778013666423
$ws.Range("H19").Value = 1.000839342965861
$ws.Range("J19").Value = 0.1267859513467897
$ws.Range("K19").Value = 0.3981692581049572
$ws.Range("M19").Value = 0.3034983265690698
$ws.Range("O19").Value = 3.986415708465813

$ws.Range("B20").Value = 0.4533652353387083
$ws.Range("C20").Value = 0.04547943327037274
$ws.Range("D20").Value = 0.6593368602205771
$ws.Range("E20").Value = 0.2584040790611439
$ws.Range("G20").Value = 0.9680470728043673
$ws.Range("H20").Value = 0.9990314015033164
$ws.Range("J20").Value = 0.126937277872706
$ws.Range("K20").Value = 0.4120687135819878
$ws.Range("M20").Value = 0.3084739930248332
$ws.Range("O20").Value = 3.981182435881436

$ws.Range("B21").Value = 0.4994131248700739
$ws.Range("C21").Value = 0.0501483170849184
$ws.Range("D21").Value = 0.6693429558087871
$ws.Range("E21").Value = 0.2611816802901643
$ws.Range("G21").Value = 0.9662298448183719
$ws.Range("H21").Value = 0.9934191834970392
$ws.Range("J21").Value = 0.1275108430352603
$ws.Range("K21").Value = 0.4587590182104009
$ws.Range("M21").Value = 0.3253166156705944
$ws.Range("O21").Value = 3.965908915027626

$ws.Range("B22").Value = 0.5295446492967528
$ws.Range("C22").Value = 0.05318803070615274
$ws.Range("D22").Value = 0.6761395718386893
$ws.Range("E22").Value = 0.2630975635133623
$ws.Range("G22").Value = 0.9655184477641541
$ws.Range("H22").Value = 0.9900954752877311
$ws.Range("J22").Value = 0.1279341820718258
$ws.Range("K22").Value = 0.489263721741338
$ws.Range("M22").Value = 0.3364164868677051
$ws.Range("O22").Value = 3.957648147073286

$ws.Range("B23").Value = 0.5134597120422768
$ws.Range("C23").Value = 0.05156676535737859
$ws.Range("D23").Value = 0.672488618926792
$ws.Range("E23").Value = 0.2620658338514659
$ws.Range("G23").Value = 0.9658545371941187
$ws.Range("H23").Value = 0.9918379255059904
$ws.Range("J23").Value = 0.1277038059160489
$ws.Range("K23").Value = 0.4729838844999961
$ws.Range("M23").Value = 0.3304838985953324
$ws.Range("O23").Value = 3.961899568898787

$ws.Range("B24").Value = 0.4526132617857286
$ws.Range("C24").Value = 0.04540292999843132
$ws.Range("D24").Value = 0.6591776664078282
$ws.Range("E24").Value = 0.2583603815173134
$ws.Range("G24").Value = 0.9680848107510514
$ws.Range("H24").Value = 0.999128945371865
$ws.Range("J24").Value = 0.1269287225726501
$ws.Range("K24").Value = 0.4113054525196844
$ws.Range("M24").Value = 0.3082002802278367
$ws.Range("O24").Value = 3.981461133605478

$ws.Range("B25").Value = 0.3871544895141596
$ws.Range("C25").Value = 0.03870603430699759
$ws.Range("D25").Value = 0.6459281906420529
$ws.Range("E25").Value = 0.2547966375913937
$ws.Range("G25").Value = 0.972533433740864
$ws.Range("H25").Value = 1.008473204339083
$ws.Range("J25").Value = 0.1263012012146589
$ws.Range("K25").Value = 0.3447490549276608
$ws.Range("M25").Value = 0.2845660137116823
$ws.Range("O25").Value = 4.009951738581805
